# "switched in depth chart projections"
# 1. Correct the spelling of "Yasmani Tomas" -> "Yasmany Tomas" (row 472, team "bears").
# 2. Insert a new player row ("Craig Gentry", team "deener") before the old row 527,
#    shifting the existing bench rows down by one.
# 3. Append 24 new bench/depth-chart players at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. fix the misspelled player name ---
$ws.Cells.Item(472, 2).Value2 = "Yasmany Tomas"

# --- 2. insert "Craig Gentry" as the new row 527 ---
$ws.Rows.Item(527).Insert()
$ws.Cells.Item(527, 1).Value2 = "deener"
$ws.Cells.Item(527, 2).Value2 = "Craig Gentry"
$ws.Cells.Item(527, 3).Value2 = 0
$ws.Cells.Item(527, 4).Value2 = "B"

# --- 3. append the new depth-chart players (rows 556-579) ---
$newPlayers = @(
    @("deener", "Anthony DeSclafini", 0, "B"),
    @("jobu", "Darin Ruf", 0, "B"),
    @("ottawa", "Jesse Chavez", 0, "B"),
    @("balco", "Kyle Parker", 0, "B"),
    @("balco", "Zach Walters", 0, "B"),
    @("ottawa", "Kevin Pillar", 0, "B"),
    @("jobu", "Dylan Axelrod", 0, "B"),
    @("deener", "Cory Luebke", 0, "B"),
    @("bellevegas", "Freddy Galvis", 0, "B"),
    @("marmaduke", "Alberto Callaspo", 0, "CI"),
    @("pasadena", "Stephen Piscotty", 0, "B"),
    @("virginia", "Chase Anderson", 0, "B"),
    @("pk dodgers", "Guillermo Heredia", 0, "B"),
    @("deano", "Justin Ruggiano", 0, "B"),
    @("bears", "Alex Jackson", 0, "B"),
    @("sturgeon", "Mark Canha", 0, "B"),
    @("dembums", "Alex Guerrero", 0, "B"),
    @("baycity", "Eduardo Escobar", 0, "B"),
    @("rippe", "Jarred Cosart", 0, "B"),
    @("rippe", "Chris Taylor", 0, "B"),
    @("d&s", "John Axford", 0, "B"),
    @("baycity", "Vance Worley", 0, "B"),
    @("dembums", "D.J. Peterson", 0, "B"),
    @("d&s", "Danny Farquhar", 0, "B")
)

$startRow = 556
for ($i = 0; $i -lt $newPlayers.Count; $i++) {
    $r = $startRow + $i
    $row = $newPlayers[$i]
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
}

# --- 4. update the view: scroll/select to match where the editor ended up ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 561
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D566").Select() | Out-Null
